$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C6").Value = "Email  is sent."
$ws.Range("B7").Value = "Email activated by admin"
$ws.Range("C7").Value = "Activates the trader and trader can login."

$ws.Activate()
$ws.Range("C7").Select()
